$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -10
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -3
$ws.Range("F8").Value = -5
$ws.Range("F10").Value = 4
$ws.Range("F12").Value = -8
$ws.Range("F14").Value = -7
$ws.Range("F16").Value = 3
